# Restructure the "PRACTICA_1_Instalacion_de_un_SOTR/" .. "El archivo
# marte_2.0_22Feb2017_src.tar.gz" region:
#   - give the "PRACTICA_1_Instalacion_de_un_SOTR/" paragraph a paragraph-mark
#     run-style (Segoe UI / color 586069) and drop the _GoBack bookmark from it
#   - insert a new paragraph "En: https://github.com/sotrteacher/sotr_201808_201812"
#   - drop the blank paragraph that used to follow
#   - re-home the _GoBack bookmark inside "El archivo ..." -> "e" | bookmark | "l archivo ..."
$d = $word.ActiveDocument

# Locate the start of the "PRACTICA_1_Instalacion_de_un_SOTR/" paragraph.
$startRange = $d.Content
$startRange.Find.Execute("PRACTICA_1_Instalacion_de_un_SOTR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startRange.Expand(4) | Out-Null

# Locate the end of the "El archivo marte_2.0_22Feb2017_src.tar.gz" paragraph.
$endRange = $d.Content
$endRange.Find.Execute("El archivo marte_2.0_22Feb2017_src.tar.gz", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRange.Expand(4) | Out-Null

# Replace that whole span (4 paragraphs) with the new 4-paragraph content.
$target = $d.Range($startRange.Start, $endRange.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="586069"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="Textoennegrita"/><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/></w:rPr><w:t>PRACTICA_1_Instalacion_de_un_SOTR</w:t></w:r><w:r><w:rPr><w:rStyle w:val="separator"/><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="586069"/></w:rPr><w:t>/</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">En: </w:t></w:r><w:r><w:t>https://github.com/sotrteacher/sotr_201808_201812</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3 Descargar las fuentes del sistema </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MaRTE</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> OS.</w:t></w:r><w:r><w:t xml:space="preserve"> Esto es, descargar </w:t></w:r></w:p><w:p><w:r><w:t>e</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>l archivo marte_2.0_22Feb2017_src.tar.gz</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml"><pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="character" w:styleId="Textoennegrita"><w:name w:val="Strong"/><w:basedOn w:val="Fuentedeprrafopredeter"/><w:uiPriority w:val="22"/><w:qFormat/><w:rPr><w:b/><w:bCs/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="separator"><w:name w:val="separator"/><w:basedOn w:val="Fuentedeprrafopredeter"/></w:style><w:style w:type="character" w:default="1" w:styleId="Fuentedeprrafopredeter"><w:name w:val="Default Paragraph Font"/></w:style></w:styles></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
